$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'29.327.86"
$ws.Range("E2").Value = "  +0.35%  "
$ws.Range("D3").Value = "'1.875.90"
$ws.Range("E3").Value = "  +0.39%  "
$ws.Range("E4").Value = "  -0.08%  "
$ws.Range("D5").Value = "'0.7119"
$ws.Range("E5").Value = "  +0.35%  "
$ws.Range("D6").Value = "'242.05"
$ws.Range("E6").Value = "  +0.36%  "
$ws.Range("D7").Value = "'1.000"
$ws.Range("E7").Value = "  -0.05%  "
$ws.Range("D8").Value = "'0.07875"
$ws.Range("E8").Value = "  +2.69%  "
$ws.Range("D9").Value = "'0.3119"
$ws.Range("E9").Value = "  +0.60%  "
$ws.Range("D10").Value = "'25.31"
$ws.Range("E10").Value = "  +1.62%  "
$ws.Range("D11").Value = "'0.08388"
$ws.Range("E11").Value = "  +0.45%  "
$ws.Range("D12").Value = "'1.867.39"
$ws.Range("E12").Value = "  -1.02%  "
$ws.Range("D13").Value = "'5.252"
$ws.Range("E13").Value = "  +0.95%  "
$ws.Range("D14").Value = "'0.7191"
$ws.Range("E14").Value = "  +1.52%  "
$ws.Range("D15").Value = "'91.40"
$ws.Range("E15").Value = "  +0.48%  "
$ws.Range("D16").Value = "'6.177"
$ws.Range("E16").Value = "  +4.32%  "
$ws.Range("D17").Value = "'0.000008370"
$ws.Range("E17").Value = "  +1.65%  "
$ws.Range("D18").Value = "'29.334.83"
$ws.Range("E18").Value = "  +0.28%  "
$ws.Range("D19").Value = "'240.92"
$ws.Range("E19").Value = "  -0.41%  "
$ws.Range("B20").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C20").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D20").Value = "'2.134.20"
$ws.Range("E20").Value = "  +0.17%  "
$ws.Range("B21").Value = "Avalanche"
$ws.Range("C21").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range("D21").Value = "'13.24"
$ws.Range("E21").Value = "  +0.83%  "
$ws.Range("D22").Value = "'1.003"
$ws.Range("E22").Value = "  +0.30%  "
$ws.Range("D23").Value = "'7.797"
$ws.Range("E23").Value = "  -0.16%  "
$ws.Range("E24").Value = "  -0.10%  "
$ws.Range("D25").Value = "'0.1601"
$ws.Range("D26").Value = "'162.96"
$ws.Range("E26").Value = "  -0.14%  "
$ws.Range("D27").Value = "'9.062"
$ws.Range("E27").Value = "  +0.86%  "
$ws.Range("E28").Value = "  +0.39%  "
$ws.Range("E29").Value = "  +0.29%  "
$ws.Range("E30").Value = "  +0.49%  "
$ws.Range("E31").Value = "  +0.59%  "
$ws.Range("D32").Value = "'1.232"
$ws.Range("E32").Value = "  -3.56%  "
$ws.Range("D33").Value = "'0.05356"
$ws.Range("E33").Value = "  +2.36%  "
$ws.Range("D34").Value = "'1.950"
$ws.Range("E34").Value = "  +1.57%  "
$ws.Range("D35").Value = "'1.179"
$ws.Range("E35").Value = "  +0.95%  "
$ws.Range("D36").Value = "'0.7477"
$ws.Range("E36").Value = "  -0.41%  "
$ws.Range("E37").Value = "  +0.18%  "
$ws.Range("D38").Value = "'1.302.14"
$ws.Range("E38").Value = "  +13.09%  "
$ws.Range("D39").Value = "'0.01884"
$ws.Range("E39").Value = "  +1.74%  "
$ws.Range("D40").Value = "'2.732"
$ws.Range("E40").Value = "  +0.72%  "
$ws.Range("D41").Value = "'6.524"
$ws.Range("E41").Value = "  +2.68%  "
$ws.Range("D42").Value = "'110.83"
$ws.Range("E42").Value = "  +6.07%  "
$ws.Range("D43").Value = "'0.8922"
$ws.Range("E43").Value = "  +0.84%  "
$ws.Range("D44").Value = "'72.84"
$ws.Range("E44").Value = "  -0.06%  "
$ws.Range("E45").Value = "  +13.99%  "
$ws.Range("D46").Value = "'0.9999"
$ws.Range("E46").Value = "  -0.02%  "
$ws.Range("D47").Value = "'2.016.11"
$ws.Range("E47").Value = "  -0.54%  "
$ws.Range("D48").Value = "'1.805"
$ws.Range("E48").Value = "  +0.93%  "
$ws.Range("E49").Value = "  +0.05%  "
$ws.Range("D50").Value = "'9.462"
$ws.Range("E50").Value = "  +1.19%  "
$ws.Range("D51").Value = "'0.4364"
$ws.Range("E51").Value = "  +1.96%  "
